# V1000_LV_HD_MODS.xlsx - build the verification parameter list.
#
# The "C6-01 / DUTY CYCLE" parameter is being dropped from the modified-
# parameter list (it isn't one of the values the programmer needs to verify
# against the VFD), so its row is deleted outright and every row below it
# shifts up to close the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (PARAMETER NUMBER "C6-01", "DUTY CYCLE", value 0) entirely -
# rows 4:19 shift up to become rows 3:18.
$ws.Rows(3).Delete()

# Column A is widest now that "CARRIER FREQUENCY UPPER LIMIT" (etc.) is the
# longest label remaining, so best-fit it to the new contents.
$ws.Columns("A:A").AutoFit()

# Leave the selection where the author was working.
$ws.Range("B10").Select()
